$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: "...thời gian :3-4 tiếng" -> "...thời gian :1 ngày"
#   Run("3-4") + Run(" tiếng")  =>  Run("1 ngày")
#   (the preceding Run(" :") must stay a separate run, untouched)
# ------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.Execute("3-4 tiếng") | Out-Null
$start1 = $find1.Start
$end1   = $find1.End

$editRng1 = $d.Range($start1, $end1)
$editRng1.Text = "1 ngày"

# Re-select the just-inserted text and nudge its formatting (off/on) so the
# engine gives it its own run instead of silently coalescing it into the
# preceding ": " run, which happens to share identical run formatting.
$newText1 = "1 ngày"
$splitRng1 = $d.Range($start1, $start1 + $newText1.Length)
$splitRng1.Font.Italic = $false
$splitRng1.Font.Italic = $true

# ------------------------------------------------------------------
# Change 2: "...thời gian: 1 ngày" -> "...thời gian: 5 ngày"
#   Run(": 1 ngày")  =>  Run(": ") + Run("5") + Run(" ngày")
# ------------------------------------------------------------------
$find2 = $d.Content
$find2.Find.Execute(": 1 ngày") | Out-Null
$start2 = $find2.Start

# Replace the digit itself.
$digitRng = $d.Range($start2 + 2, $start2 + 3)
$digitRng.Text = "5"

# Force the "5" and the trailing " ngày" into their own runs (same trick as
# above) so the three pieces stay distinct instead of re-merging.
$splitA = $d.Range($start2 + 2, $start2 + 3)   # "5"
$splitA.Font.Italic = $false
$splitA.Font.Italic = $true

$splitB = $d.Range($start2 + 3, $start2 + 8)   # " ngày"
$splitB.Font.Italic = $false
$splitB.Font.Italic = $true
